$d = $word.ActiveDocument

# 1) Update "Abdon Morales" to include student id, in the team members list.
$d.Content.Find.Execute("Abdon Morales", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Abdon Morales (am226923)", 2)

# 2) Append a new paragraph at the end of the document with the same text,
#    mirroring the addition after "Lauren N Parker".
$d.Content.InsertAfter("`r`nAbdon Morales (am226923)")
